# chore: update Sheets via scheduled runner
# Refresh computed price/profit columns (H-N) across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets to reflect the latest market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 267310.25
$ws.Range("J17").Value = 267310.25
$ws.Range("L17").Value = 801930.75
$ws.Range("N17").Value = -802266.75
$ws.Range("H112").Value = 1419.5555
$ws.Range("I112").Value = 933.2
$ws.Range("J112").Value = 1606.6154
$ws.Range("K112").Value = 2799.6
$ws.Range("L112").Value = 4819.8462
$ws.Range("M112").Value = -1691.6
$ws.Range("N112").Value = -7035.8462
$ws.Range("H132").Value = 1381.6333
$ws.Range("I132").Value = 1321.5714
$ws.Range("K132").Value = 3964.7142
$ws.Range("M132").Value = -1434.7142
$ws.Range("H134").Value = 125916.71
$ws.Range("J134").Value = 116986.5
$ws.Range("L134").Value = 116986.5
$ws.Range("N134").Value = -127126.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 5016.276
$ws.Range("I23").Value = 3945.4285
$ws.Range("K23").Value = 3945.4285
$ws.Range("M23").Value = -3686.4285
$ws.Range("H32").Value = 16427.354
$ws.Range("I32").Value = 3329.1125
$ws.Range("K32").Value = 3329.1125
$ws.Range("M32").Value = -3042.1125
$ws.Range("H45").Value = 317599.25
$ws.Range("I45").Value = 440692.44
$ws.Range("J45").Value = 3027.7778
$ws.Range("K45").Value = 440692.44
$ws.Range("L45").Value = 3027.7778
$ws.Range("M45").Value = -440315.44
$ws.Range("N45").Value = -3781.7778
$ws.Range("H61").Value = 1334.625
$ws.Range("I61").Value = 1345.3871
$ws.Range("K61").Value = 1345.3871
$ws.Range("M61").Value = -1133.3871
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9470
$ws.Range("H136").Value = 1334.625
$ws.Range("I136").Value = 1345.3871
$ws.Range("K136").Value = 4036.1613
$ws.Range("M136").Value = -1486.1613
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8114.1035
$ws.Range("I20").Value = 6886.85
$ws.Range("J20").Value = 10841.333
$ws.Range("K20").Value = 6886.85
$ws.Range("L20").Value = 10841.333
$ws.Range("M20").Value = -6639.85
$ws.Range("N20").Value = -11335.333
$ws.Range("H105").Value = 4951.5386
$ws.Range("I105").Value = 5670.1816
$ws.Range("K105").Value = 5670.1816
$ws.Range("M105").Value = -3923.1816
$ws.Range("H134").Value = 1203.0667
$ws.Range("I134").Value = 1203.0667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3609.2001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1074.2001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2641.0715
$ws.Range("I16").Value = 3131.111
$ws.Range("K16").Value = 3131.111
$ws.Range("M16").Value = -2844.111
$ws.Range("H22").Value = 563.5714
$ws.Range("I22").Value = 411.1111
$ws.Range("J22").Value = 838
$ws.Range("K22").Value = 411.1111
$ws.Range("L22").Value = 838
$ws.Range("M22").Value = -61.11110000000002
$ws.Range("N22").Value = -1538
$ws.Range("H113").Value = 2641.0715
$ws.Range("I113").Value = 3131.111
$ws.Range("K113").Value = 3131.111
$ws.Range("M113").Value = -961.1109999999999
$ws.Range("H132").Value = 4145.5
$ws.Range("I132").Value = 3935.675
$ws.Range("J132").Value = 4984.8
$ws.Range("K132").Value = 11807.025
$ws.Range("L132").Value = 14954.4
$ws.Range("M132").Value = -9277.025000000001
$ws.Range("N132").Value = -20014.4
$ws.Range("H134").Value = 2294.1228
$ws.Range("I134").Value = 2167.9556
$ws.Range("K134").Value = 6503.8668
$ws.Range("M134").Value = -3968.8668
$ws.Range("H141").Value = 172133
$ws.Range("J141").Value = 172133
$ws.Range("L141").Value = 172133
$ws.Range("N141").Value = -182493

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 471.83334
$ws.Range("I14").Value = 471.83334
$ws.Range("K14").Value = 1415.50002
$ws.Range("M14").Value = -1242.50002
$ws.Range("H107").Value = 636.4681
$ws.Range("J107").Value = 641.9
$ws.Range("L107").Value = 1925.7
$ws.Range("N107").Value = -5765.7
$ws.Range("H121").Value = 13973922
$ws.Range("I121").Value = 41792030
$ws.Range("K121").Value = 125376090
$ws.Range("M121").Value = -125374780
$ws.Range("H132").Value = 2266.5833
$ws.Range("I132").Value = 1747.25
$ws.Range("K132").Value = 15725.25
$ws.Range("M132").Value = -13195.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5209
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 5209
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H102").Value = 3241.0454
$ws.Range("I102").Value = 3241.0454
$ws.Range("K102").Value = 3241.0454
$ws.Range("M102").Value = -1619.0454
$ws.Range("H109").Value = 65499.5
$ws.Range("J109").Value = 65499.5
$ws.Range("L109").Value = 65499.5
$ws.Range("N109").Value = -67579.5
$ws.Range("H126").Value = 2432.5
$ws.Range("I126").Value = 2300.6667
$ws.Range("J126").Value = 2828
$ws.Range("K126").Value = 6902.000100000001
$ws.Range("L126").Value = 8484
$ws.Range("M126").Value = -4432.000100000001
$ws.Range("N126").Value = -13424
$ws.Range("H132").Value = 2826.2856
$ws.Range("I132").Value = 2861.6667
$ws.Range("J132").Value = 2762.6
$ws.Range("K132").Value = 8585.000100000001
$ws.Range("L132").Value = 8287.8
$ws.Range("M132").Value = -6055.000100000001
$ws.Range("N132").Value = -13347.8
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8324
$ws.Range("I7").Value = 10813.0625
$ws.Range("K7").Value = 10813.0625
$ws.Range("M7").Value = -10701.0625
$ws.Range("H46").Value = 141341
$ws.Range("I46").Value = 141341
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 141341
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -141153
$ws.Range("H87").Value = 25189
$ws.Range("J87").Value = 25189
$ws.Range("L87").Value = 25189
$ws.Range("N87").Value = -27435
$ws.Range("H90").Value = 25189
$ws.Range("J90").Value = 25189
$ws.Range("L90").Value = 75567
$ws.Range("N90").Value = -86799
$ws.Range("H100").Value = 25623.562
$ws.Range("I100").Value = 4960.4443
$ws.Range("K100").Value = 4960.4443
$ws.Range("M100").Value = -4419.4443
$ws.Range("H121").Value = 109999
$ws.Range("J121").Value = 109999
$ws.Range("L121").Value = 109999
$ws.Range("N121").Value = -113493
$ws.Range("H126").Value = 8324
$ws.Range("I126").Value = 10813.0625
$ws.Range("K126").Value = 32439.1875
$ws.Range("M126").Value = -29969.1875
$ws.Range("H132").Value = 2772.7314
$ws.Range("I132").Value = 2360.848
$ws.Range("K132").Value = 7082.544
$ws.Range("M132").Value = -4552.544
$ws.Range("H136").Value = 3047.0173
$ws.Range("I136").Value = 2622.6829
$ws.Range("K136").Value = 7868.048699999999
$ws.Range("M136").Value = -5318.048699999999
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11949.8
$ws.Range("I81").Value = 9388.777
$ws.Range("K81").Value = 18777.554
$ws.Range("M81").Value = -17716.554
$ws.Range("H84").Value = 11949.8
$ws.Range("I84").Value = 9388.777
$ws.Range("K84").Value = 93887.77
$ws.Range("M84").Value = -88583.77
$ws.Range("H132").Value = 3846.5
$ws.Range("I132").Value = 5482.3
$ws.Range("J132").Value = 1509.6428
$ws.Range("K132").Value = 16446.9
$ws.Range("L132").Value = 4528.928400000001
$ws.Range("M132").Value = -13916.9
$ws.Range("N132").Value = -9588.9284
$ws.Range("H136").Value = 594.8095
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
